$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "last sync" timestamp (Excel serial date number) shared by both rows.
$newSync = 45306.895564331055

# Update the sync timestamps for the two data rows.
$ws.Range("E2").Value = $newSync
$ws.Range("E3").Value = $newSync

# Remove the "Last Sync" header cell (and its now-unused shared string).
$ws.Range("E1").Clear()

# Move the active selection to E1 (mirrors the recorded selection change).
$ws.Range("E1").Select()
